$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.821.39"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "2.301.96"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "2.298.80"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "59.744.23"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "2.711.56"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "2.301.89"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "309.73"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("E22").Value = "  -8.29%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.12"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "0.0₃0711"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.66"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.04"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "310.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "136.02"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.41%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0934"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.563"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0488"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0224"
$ws.Range("E48").Value = "  +8.75%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0211"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.98"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.21%  "
